# Auto-generated edit script: updates crypto price/volume table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay stored as text (match source formatting)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply cell value updates
$ws.Range("D2").Value = "63.827.37"
$ws.Range("E2").Value = "  +0.65%  "
$ws.Range("D3").Value = "3.323.28"
$ws.Range("E3").Value = "  +1.75%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "604.99"
$ws.Range("E5").Value = "  +1.66%  "
$ws.Range("D6").Value = "143.06"
$ws.Range("E6").Value = "  +0.99%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "3.321.57"
$ws.Range("E8").Value = "  +1.84%  "
$ws.Range("D9").Value = "0.519"
$ws.Range("E9").Value = "  -0.31%  "
$ws.Range("D10").Value = "0.150"
$ws.Range("E10").Value = "  +1.07%  "
$ws.Range("E11").Value = "  +3.41%  "
$ws.Range("D12").Value = "0.470"
$ws.Range("E12").Value = "  +0.78%  "
$ws.Range("E13").Value = "  -0.40%  "
$ws.Range("D14").Value = "35.02"
$ws.Range("E14").Value = "  +0.54%  "
$ws.Range("D15").Value = "3.872.29"
$ws.Range("E15").Value = "  +1.75%  "
$ws.Range("D16").Value = "0.120"
$ws.Range("E16").Value = "  +0.36%  "
$ws.Range("D17").Value = "3.324.15"
$ws.Range("E17").Value = "  +1.74%  "
$ws.Range("D18").Value = "63.905.76"
$ws.Range("D19").Value = "6.85"
$ws.Range("E19").Value = "  +0.56%  "
$ws.Range("D20").Value = "480.73"
$ws.Range("E20").Value = "  +0.37%  "
$ws.Range("D21").Value = "14.14"
$ws.Range("E21").Value = "  -0.85%  "
$ws.Range("D22").Value = "0.737"
$ws.Range("E22").Value = "  +0.75%  "
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("D24").Value = "13.84"
$ws.Range("E24").Value = "  +3.01%  "
$ws.Range("D25").Value = "84.85"
$ws.Range("E25").Value = "  +0.60%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("E27").Value = "  +1.25%  "
$ws.Range("B28").Value = "FirstDigitalUSD"
$ws.Range("C28").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").Value = "8.27"
$ws.Range("E29").Value = "  +1.48%  "
$ws.Range("E30").Value = "  -2.39%  "
$ws.Range("D31").Value = "2.16"
$ws.Range("E31").Value = "  +0.46%  "
$ws.Range("D32").Value = "28.92"
$ws.Range("E32").Value = "  +3.06%  "
$ws.Range("E33").Value = "  -1.76%  "
$ws.Range("E34").Value = "  -1.42%  "
$ws.Range("D35").Value = "1.11"
$ws.Range("E35").Value = "  +1.06%  "
$ws.Range("D36").Value = "6.08"
$ws.Range("E36").Value = "  +2.30%  "
$ws.Range("D37").Value = "0.0₃0752"
$ws.Range("E37").Value = "  +4.05%  "
$ws.Range("D38").Value = "52.42"
$ws.Range("E38").Value = "  -1.04%  "
$ws.Range("D39").Value = "0.0401"
$ws.Range("E39").Value = "  +1.31%  "
$ws.Range("D40").Value = "433.69"
$ws.Range("E40").Value = "  +2.45%  "
$ws.Range("D41").Value = "3.108.53"
$ws.Range("E41").Value = "  +3.67%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "0.118"
$ws.Range("E42").Value = "  +4.74%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").Value = "2.77"
$ws.Range("E43").Value = "  -0.57%  "
$ws.Range("E44").Value = "  -0.87%  "
$ws.Range("D45").Value = "0.267"
$ws.Range("E45").Value = "  -0.60%  "
$ws.Range("D46").Value = "2.25"
$ws.Range("E46").Value = "  +3.87%  "
$ws.Range("D47").Value = "36.55"
$ws.Range("E47").Value = "  +8.38%  "
$ws.Range("D48").Value = "26.49"
$ws.Range("E48").Value = "  +1.45%  "
$ws.Range("D49").Value = "0.998"
$ws.Range("E49").Value = "  -0.09%  "
$ws.Range("E50").Value = "  -0.97%  "
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").Value = "124.61"
$ws.Range("E51").Value = "  +4.42%  "
